$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 11: 2021年 data (all 43 columns A:AQ) ---
$row11 = @("2021年",23812,6032,1956,5,23128,23146,2084,8629,6579,7149,2657,319,9989,8395,1171,11205,23278,3166,17254,122,4558,3029,11049,30305,8308,2229,144,19774,13168,441517,24160,30511,7189,5665,639,31287,5646,44086,3707,8874,5640,1473)

for ($i = 0; $i -lt $row11.Length; $i++) {
    $ws.Cells.Item(11, $i + 1).Value = $row11[$i]
}

# --- Row 12: 2022年, mostly blank except the aggregate total in column AE ---
$ws.Cells.Item(12, 1).Value = "2022年"
$ws.Cells.Item(12, 31).Value = 451000   # column AE = 31st column

# --- Formatting: column A labels use the same bold/centered/bordered style as existing year rows ---
$ws.Range("A10").Copy() | Out-Null
$ws.Range("A11").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("A10").Copy() | Out-Null
$ws.Range("A12").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false
